# Changes for 57748 and 57772 user stories
# Add a new Finding row (ICW_UnbornChildDOBListedInThePast) to the "Findings" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Findings")

# Insert a new row at position 14 (pushes existing rows 14+ down by one),
# then populate it with the new finding's Name / Value / Note.
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = "ICW_UnbornChildDOBListedInThePast"
$ws.Range("B14").Value = "There is a Unborn Child listed on the ICW but has a date of birth listed in the past."
$ws.Range("C14").Value = "Newly added by Raluca"

# Reflect the author's final cursor/selection position on the sheet.
$ws.Activate()
$ws.Range("B15").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
